$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.525.84'
$ws.Range('E2').Value = '  -0.69%  '
$ws.Range('D3').Value = '1.851.70'
$ws.Range('E3').Value = '  -0.29%  '
$ws.Range('D4').Value = '0.9986'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '241.91'
$ws.Range('E5').Value = '  -1.12%  '
$ws.Range('D6').Value = '0.6303'
$ws.Range('D7').Value = '0.9997'
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '47.94'
$ws.Range('E8').Value = '  +1.29%  '
$ws.Range('D9').Value = '0.07576'
$ws.Range('E9').Value = '  +0.79%  '
$ws.Range('D10').Value = '0.2979'
$ws.Range('E10').Value = '  -0.14%  '
$ws.Range('D11').Value = '24.38'
$ws.Range('E11').Value = '  -0.30%  '
$ws.Range('D12').Value = '0.07672'
$ws.Range('E12').Value = '  +0.10%  '
$ws.Range('D13').Value = '1.891.67'
$ws.Range('E13').Value = '  +1.17%  '
$ws.Range('D14').Value = '5.018'
$ws.Range('E14').Value = '  -0.39%  '
$ws.Range('D15').Value = '0.6854'
$ws.Range('E15').Value = '  -0.85%  '
$ws.Range('D16').Value = '83.83'
$ws.Range('E16').Value = '  +0.02%  '
$ws.Range('D17').Value = '0.000009810'
$ws.Range('E17').Value = '  -0.18%  '
$ws.Range('D18').Value = '2.136.14'
$ws.Range('E18').Value = '  +0.84%  '
$ws.Range('D19').Value = '6.220'
$ws.Range('E19').Value = '  +2.14%  '
$ws.Range('D20').Value = '29.566.66'
$ws.Range('E20').Value = '  -0.58%  '
$ws.Range('D21').Value = '234.55'
$ws.Range('E21').Value = '  -0.52%  '
$ws.Range('E22').Value = '  -1.42%  '
$ws.Range('D23').Value = '0.9998'
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('D24').Value = '7.625'
$ws.Range('E24').Value = '  +1.14%  '
$ws.Range('D25').Value = '0.9998'
$ws.Range('E25').Value = '  -0.07%  '
$ws.Range('D26').Value = '155.64'
$ws.Range('E26').Value = '  -2.07%  '
$ws.Range('D27').Value = '0.1391'
$ws.Range('E27').Value = '  -2.10%  '
$ws.Range('D28').Value = '8.440'
$ws.Range('E28').Value = '  -1.13%  '
$ws.Range('D29').Value = '17.73'
$ws.Range('E29').Value = '  -1.00%  '
$ws.Range('D30').Value = '1.480'
$ws.Range('E30').Value = '  -0.80%  '
$ws.Range('D31').Value = '0.05839'
$ws.Range('E31').Value = '  -5.81%  '
$ws.Range('D32').Value = '1.261'
$ws.Range('E32').Value = '  -1.61%  '
$ws.Range('D33').Value = '4.111'
$ws.Range('E33').Value = '  -1.24%  '
$ws.Range('D34').Value = '4.043'
$ws.Range('E34').Value = '  -1.35%  '
$ws.Range('D35').Value = '1.916'
$ws.Range('E35').Value = '  +1.10%  '
$ws.Range('E36').Value = '  +0.04%  '
$ws.Range('D37').Value = '0.7174'
$ws.Range('E37').Value = '  -1.40%  '
$ws.Range('D38').Value = '2.587'
$ws.Range('E38').Value = '  -0.62%  '
$ws.Range('D40').Value = '1.234.70'
$ws.Range('E40').Value = '  +2.60%  '
$ws.Range('E41').Value = '  -0.31%  '
$ws.Range('D42').Value = '0.9137'
$ws.Range('E42').Value = '  -1.19%  '
$ws.Range('D43').Value = '6.130'
$ws.Range('E43').Value = '  -1.62%  '
$ws.Range('D44').Value = '2.037.85'
$ws.Range('E44').Value = '  +0.48%  '
$ws.Range('E45').Value = '  -0.03%  '
$ws.Range('D46').Value = '101.91'
$ws.Range('E46').Value = '  -0.06%  '
$ws.Range('E47').Value = '  +1.34%  '
$ws.Range('D48').Value = '7.301'
$ws.Range('E48').Value = '  +9.25%  '
$ws.Range('D49').Value = '9.178'
$ws.Range('E49').Value = '  -0.06%  '
$ws.Range('E50').Value = '  -1.29%  '
$ws.Range('D51').Value = '0.4032'
$ws.Range('E51').Value = '  -0.69%  '
